$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new data rows right after the first worker row (old row 16),
# pushing the old second worker row (17) down to 19 and the signature block
# (old rows 22/23) down to 24/25.
# ---------------------------------------------------------------------------
$ws.Rows("17:18").Insert()

# Copy the formatting (borders/styles) of row 16 down into the two freshly
# inserted rows so they look like ordinary (non-final) data rows.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J18").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# Header figures
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 132411
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 4

# ---------------------------------------------------------------------------
# Data rows (2:10 columns B..J) -> Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico
# ---------------------------------------------------------------------------
# Row 16: YUDIS MARIA VILLERO TOVAR, periodo 2307
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "45553717"
$ws.Range("D16").Value = "YUDIS MARIA VILLERO TOVAR"
$ws.Range("E16").Value = "2307"
$ws.Range("F16").Value = 46400
$ws.Range("G16").Value = 1160000

# Row 17 (new): YUDIS MARIA VILLERO TOVAR, periodo 2306
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45553717"
$ws.Range("D17").Value = "YUDIS MARIA VILLERO TOVAR"
$ws.Range("E17").Value = "2306"
$ws.Range("F17").Value = 46400
$ws.Range("G17").Value = 1160000

# Row 18 (new): ROMAN MEZA, periodo 2111 (formerly row 16's content)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "79376090"
$ws.Range("D18").Value = "ROMAN MEZA"
$ws.Range("E18").Value = "2111"
$ws.Range("F18").Value = 1211
$ws.Range("G18").Value = 908526

# Row 19: RAUL ENRIQUE VELEZ TATIS (formerly row 17, values unchanged, already shifted)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "20341526"
$ws.Range("D19").Value = "RAUL ENRIQUE VELEZ TATIS"
$ws.Range("E19").Value = "2201"
$ws.Range("F19").Value = 38400
$ws.Range("G19").Value = 1200000
